$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet grows from 42 data rows to 48; extend the formatting (border,
# wrap text, alignment) of the existing data rows down into the newly
# added rows 43:48 before writing values into them, so the new rows look
# identical to the existing ones.
$formatSource = $ws.Range("A2:C2")
$formatTarget = $ws.Range("A43:C48")
$formatSource.Copy()
$formatTarget.PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$data = @(
    @('''4812969', 'https://vnexpress.net/quy-mo-thuong-mai-dien-tu-viet-nam-uoc-dat-22-ty-usd-4812969.html', 'Quy mô thương mại điện tử Việt Nam ước đạt 22 tỷ USD'),
    @('''4812762', 'https://vnexpress.net/ong-donald-trump-hua-hen-nhung-gi-voi-kinh-te-my-4812762.html', 'Ông Donald Trump hứa hẹn những gì với kinh tế Mỹ?'),
    @('''4812976', 'https://vnexpress.net/tuyen-duong-sat-noi-voi-trung-quoc-du-kien-khoi-cong-vao-2025-4812976.html', 'Tuyến đường sắt nối với Trung Quốc dự kiến khởi công vào 2025'),
    @('''4812911', 'https://vnexpress.net/nhan-tron-giam-ve-vung-88-trieu-dong-4812911.html', 'Vàng nhẫn trơn giảm về quanh 88 triệu đồng'),
    @('''4812924', 'https://vnexpress.net/chung-khoan-tang-15-diem-4812924.html', 'Chứng khoán tăng 15 điểm'),
    @('''4812646', 'https://vnexpress.net/bitcoin-cham-71-500-usd-khi-trump-tam-dan-dau-4812646.html', 'Bitcoin lập kỷ lục 75.300 USD'),
    @('''4812599', 'https://vnexpress.net/nvidia-vuot-apple-de-thanh-cong-ty-gia-tri-nhat-the-gioi-4812599.html', 'Nvidia vượt Apple để thành công ty giá trị nhất thế giới'),
    @('''4812363', 'https://vnexpress.net/ethiopia-gap-kho-khi-cam-xe-xang-chuyen-sang-oto-dien-4812363.html', 'Ethiopia gặp khó khi cấm xe xăng chuyển sang ôtô điện'),
    @('''4812502', 'https://vnexpress.net/thu-thue-tu-chuyen-nhuong-dat-tang-manh-4812502.html', 'Thu thuế chuyển nhượng bất động sản ở TP HCM tăng gấp rưỡi'),
    @('''4812977', 'https://vnexpress.net/mobile-money-co-gan-10-trieu-tai-khoan-sau-3-nam-thi-diem-4812977.html', 'Mobile Money có gần 10 triệu tài khoản sau 3 năm thí điểm'),
    @('''4812645', 'https://vnexpress.net/kiem-soat-quyen-luc-de-giam-ganh-nang-cho-nguoi-duyet-du-an-4812645.html', '''''Kiểm soát quyền lực để giảm gánh nặng cho người duyệt dự án'''),
    @('''4812590', 'https://vnexpress.net/gia-vang-chung-khoan-my-cung-tang-trong-ngay-bau-cu-4812590.html', 'Giá vàng, chứng khoán Mỹ cùng tăng trong ngày bầu cử'),
    @('''4811187', 'https://vnexpress.net/temu-phai-nop-thue-the-nao-o-viet-nam-4811187.html', 'Temu phải nộp thuế thế nào ở Việt Nam'),
    @('''4812563', 'https://vnexpress.net/viet-nam-co-hon-9-trieu-tai-khoan-chung-khoan-4812563.html', 'Việt Nam có hơn 9 triệu tài khoản chứng khoán'),
    @('''4812805', 'https://vnexpress.net/nha-may-lego-1-3-ty-usd-chay-thu-nghiem-4812805.html', 'Nhà máy Lego 1,3 tỷ USD chạy thử nghiệm'),
    @('''4812719', 'https://vnexpress.net/vinhomes-va-vinfast-la-thuong-hieu-quoc-gia-4812719.html', 'Vinhomes và VinFast là thương hiệu quốc gia'),
    @('''4812510', 'https://vnexpress.net/vib-nhan-ky-luc-ve-tinh-nang-ca-nhan-hoa-thiet-ke-the-tin-dung-4812510.html', 'VIB nhận kỷ lục về tính năng cá nhân hóa thiết kế thẻ tín dụng'),
    @('''4812450', 'https://vnexpress.net/hai-phong-giai-ngan-gan-9-500-ty-von-dau-tu-cong-4812450.html', 'Hải Phòng giải ngân gần 9.500 tỷ vốn đầu tư công'),
    @('''4812537', 'https://vnexpress.net/keo-dai-thoi-gian-thi-diem-cho-nguoi-viet-vao-choi-casino-4812537.html', 'Kéo dài thời gian thí điểm cho người Việt vào chơi casino'),
    @('''4812530', 'https://vnexpress.net/doanh-nghiep-muon-mua-du-an-cua-ong-dung-lo-voi-bao-lo-4812530.html', 'Doanh nghiệp muốn mua dự án của ông Dũng ''Lò Vôi'' báo lỗ'),
    @('''4812449', 'https://vnexpress.net/chung-khoan-hom-nay-5-11-vn-index-co-phien-ru-ngu-nha-dau-tu-4812449.html', 'Thanh khoản chứng khoán giảm mạnh'),
    @('''4812386', 'https://vnexpress.net/cong-nhan-boeing-dung-dinh-cong-4812386.html', 'Công nhân Boeing dừng đình công'),
    @('''4812352', 'https://vnexpress.net/pho-thu-tuong-se-dung-ai-kiem-soat-doanh-thu-tren-san-thuong-mai-dien-tu-4812352.html', 'Phó thủ tướng: Sẽ dùng AI kiểm soát doanh thu trên sàn thương mại điện tử'),
    @('''4812272', 'https://vnexpress.net/viet-nam-chi-gan-1-2-ty-usd-nhap-khau-gao-4812272.html', 'Việt Nam chi gần 1,2 tỷ USD nhập khẩu gạo'),
    @('''4811603', 'https://vnexpress.net/xuat-khau-sang-my-the-nao-khi-trump-hoac-harris-thang-cu-4811603.html', 'Xuất khẩu sang Mỹ thế nào khi Trump hoặc Harris thắng cử'),
    @('''4811459', 'https://vnexpress.net/kich-ban-gia-vang-va-chung-khoan-sau-ky-bau-cu-my-4811459.html', 'Kịch bản giá vàng và chứng khoán sau kỳ bầu cử Mỹ'),
    @('''4812506', 'https://vnexpress.net/hang-khong-dang-phuc-hoi-4812506.html', 'Hàng không đang phục hồi'),
    @('''4812302', 'https://vnexpress.net/nganh-ban-le-my-canh-bao-thue-nhap-khau-cua-trump-4812302.html', 'Ngành bán lẻ Mỹ cảnh báo thuế nhập khẩu của Trump'),
    @('''4811919', 'https://vnexpress.net/thiet-bi-dien-mat-troi-trung-quoc-choi-meo-von-chuot-voi-my-4811919.html', 'Thiết bị điện mặt trời Trung Quốc chơi ''mèo vờn chuột'' với Mỹ'),
    @('''4812199', 'https://vnexpress.net/co-phieu-trump-media-tang-vot-truoc-ngay-bau-cu-my-4812199.html', 'Cổ phiếu Trump Media tăng vọt trước ngày bầu cử Mỹ'),
    @('''4811974', 'https://vnexpress.net/ong-trump-khien-dong-von-tram-ty-usd-vao-my-thap-thom-4811974.html', 'Ông Trump khiến dòng vốn trăm tỷ USD vào Mỹ thấp thỏm'),
    @('''4811962', 'https://vnexpress.net/gia-ca-phe-quay-dau-giam-4811962.html', 'Giá cà phê quay đầu giảm'),
    @('''4812439', 'https://vnexpress.net/loat-hoat-dong-ky-niem-sinh-nhat-20-nam-cua-vincom-4812439.html', 'Loạt hoạt động kỷ niệm sinh nhật 20 năm của Vincom'),
    @('''4812448', 'https://vnexpress.net/yeu-to-giup-bidv-dat-thuong-hieu-quoc-gia-4812448.html', 'Yếu tố giúp BIDV đạt thương hiệu quốc gia'),
    @('''4812496', 'https://vnexpress.net/5-san-pham-cua-mobifone-la-thuong-hieu-quoc-gia-viet-nam-2024-4812496.html', '5 sản phẩm của MobiFone là Thương hiệu quốc gia Việt Nam 2024'),
    @('''4812443', 'https://vnexpress.net/agribank-dat-thuong-hieu-quoc-gia-nho-hoat-dong-hieu-qua-4812443.html', 'Agribank đạt thương hiệu quốc gia nhờ hoạt động hiệu quả'),
    @('''4812327', 'https://vnexpress.net/chien-luoc-xanh-giup-hdbank-dat-thuong-hieu-quoc-gia-4812327.html', 'Chiến lược xanh giúp HDBank đạt thương hiệu quốc gia'),
    @('''4812100', 'https://vnexpress.net/de-xuat-giam-tiep-thue-moi-truong-voi-xang-dau-het-nam-2025-4812100.html', 'Đề xuất giảm tiếp thuế môi trường với xăng dầu hết năm 2025'),
    @('''4812088', 'https://vnexpress.net/evn-de-xuat-thi-diem-gia-dien-hai-thanh-phan-4812088.html', 'EVN đề xuất thí điểm giá điện hai thành phần'),
    @('''4811845', 'https://vnexpress.net/opec-lai-hoan-ke-hoach-tang-bom-dau-4811845.html', 'OPEC+ lại hoãn kế hoạch tăng bơm dầu'),
    @('''4812022', 'https://vnexpress.net/doanh-nghiep-tu-nhan-can-co-che-chu-khong-phai-tien-4812022.html', '''''Doanh nghiệp tư nhân cần cơ chế chứ không phải tiền'''),
    @('''4812044', 'https://vnexpress.net/chung-khoan-do-lua-vi-co-phieu-ngan-hang-4812044.html', 'Chứng khoán đỏ lửa vì cổ phiếu ngân hàng'),
    @('''4811999', 'https://vnexpress.net/vang-mieng-giam-nua-trieu-dong-4811999.html', 'Giá vàng miếng giảm, USD lên kịch trần'),
    @('''4811981', 'https://vnexpress.net/ty-phu-thai-nhan-hon-12-000-ty-dong-co-tuc-sau-7-nam-thau-tom-sabeco-4811981.html', 'Tỷ phú Thái nhận hơn 12.000 tỷ đồng cổ tức sau 7 năm thâu tóm Sabeco'),
    @('''4811774', 'https://vnexpress.net/nganh-dien-gio-my-lo-trump-tai-dac-cu-4811774.html', 'Ngành điện gió Mỹ lo Trump tái đắc cử'),
    @('''4811849', 'https://vnexpress.net/lang-phi-hang-tram-nghin-ty-dong-khi-du-an-dau-tu-roi-dap-chieu-4811849.html', 'Lãng phí hàng trăm nghìn tỷ đồng khi dự án đầu tư rồi ''đắp chiếu'''),
    @('''4811828', 'https://vnexpress.net/mot-doanh-nghiep-may-lon-khong-co-don-hang-suot-18-thang-4811828.html', 'Một doanh nghiệp may lớn không có đơn hàng suốt 18 tháng'),
)

$row = 2
foreach ($item in $data) {
    # Values below are pre-escaped: column A gets a leading apostrophe so
    # Excel stores the numeric-looking article ID as text (matching the
    # workbook's original inlineStr storage); any title/URL that itself
    # starts with a literal apostrophe has an extra one prepended, since
    # Excel strips exactly one leading apostrophe as its quote-prefix marker.
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]
    $row++
}

Write-Host "Updated rows 2 to $($row - 1)"